$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2197.1304
$ws.Range("I15").Value = 2197.1304
$ws.Range("K15").Value = 6591.3912
$ws.Range("M15").Value = -6422.3912
$ws.Range("H28").Value = 937.425
$ws.Range("I28").Value = 812.5862
$ws.Range("K28").Value = 812.5862
$ws.Range("M28").Value = -327.5862
$ws.Range("H34").Value = 2861
$ws.Range("I34").Value = 2861
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2861
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2658
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 2861
$ws.Range("I36").Value = 2861
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2861
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2146
$ws.Range("N36").ClearContents()
$ws.Range("H47").Value = 7531
$ws.Range("I47").Value = 10067
$ws.Range("K47").Value = 10067
$ws.Range("M47").Value = -9095
$ws.Range("H54").Value = 14038
$ws.Range("I54").Value = 8076
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 8076
$ws.Range("L54").Value = 20000
$ws.Range("M54").Value = -7590
$ws.Range("N54").Value = -20972
$ws.Range("H113").Value = 5743.625
$ws.Range("I113").Value = 3669.875
$ws.Range("J113").Value = 7817.375
$ws.Range("K113").Value = 3669.875
$ws.Range("L113").Value = 7817.375
$ws.Range("M113").Value = -415.875
$ws.Range("N113").Value = -14325.375
$ws.Range("H132").Value = 27526.762
$ws.Range("I132").Value = 1753.9375
$ws.Range("K132").Value = 5261.8125
$ws.Range("M132").Value = -2731.8125
$ws.Range("H135").Value = 2949.3809
$ws.Range("I135").Value = 3183.7896
$ws.Range("K135").Value = 28654.1064
$ws.Range("M135").Value = -26119.1064
$ws.Range("H138").Value = 1340.6666
$ws.Range("I138").Value = 954.73914
$ws.Range("J138").Value = 3559.75
$ws.Range("K138").Value = 2864.21742
$ws.Range("L138").Value = 10679.25
$ws.Range("M138").Value = 2275.78258
$ws.Range("N138").Value = -20959.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 750
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 750
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 750
$ws.Range("N10").Value = -1090
$ws.Range("M10").ClearContents()
$ws.Range("H32").Value = 8147.5137
$ws.Range("I32").Value = 7207.1943
$ws.Range("K32").Value = 7207.1943
$ws.Range("M32").Value = -6920.1943
$ws.Range("H45").Value = 2532.4443
$ws.Range("I45").Value = 1711.3334
$ws.Range("K45").Value = 1711.3334
$ws.Range("M45").Value = -1334.3334
$ws.Range("H74").Value = 750.2692
$ws.Range("I74").Value = 735.087
$ws.Range("K74").Value = 735.087
$ws.Range("M74").Value = 138.913
$ws.Range("H77").Value = 750.2692
$ws.Range("I77").Value = 735.087
$ws.Range("K77").Value = 3675.435
$ws.Range("M77").Value = 692.5650000000001
$ws.Range("H97").Value = 3704524
$ws.Range("I97").Value = 896.1818
$ws.Range("K97").Value = 896.1818
$ws.Range("M97").Value = -400.1818
$ws.Range("H122").Value = 3712.1765
$ws.Range("I122").Value = 2118.0435
$ws.Range("K122").Value = 6354.130500000001
$ws.Range("M122").Value = -3904.130500000001
$ws.Range("H132").Value = 2954.2856
$ws.Range("I132").Value = 1913.1578
$ws.Range("K132").Value = 5739.4734
$ws.Range("M132").Value = -3209.4734

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3279.7878
$ws.Range("I134").Value = 2215
$ws.Range("J134").Value = 5728.8
$ws.Range("K134").Value = 6645
$ws.Range("L134").Value = 17186.4
$ws.Range("M134").Value = -4110
$ws.Range("N134").Value = -22256.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 2562.5
$ws.Range("J23").Value = 2562.5
$ws.Range("L23").Value = 2562.5
$ws.Range("N23").Value = -3042.5
$ws.Range("H27").Value = 2562.5
$ws.Range("J27").Value = 2562.5
$ws.Range("L27").Value = 2562.5
$ws.Range("N27").Value = -2946.5
$ws.Range("H31").Value = 4443.5557
$ws.Range("I31").Value = 1981.6666
$ws.Range("J31").Value = 6905.4443
$ws.Range("K31").Value = 1981.6666
$ws.Range("L31").Value = 6905.4443
$ws.Range("M31").Value = -1686.6666
$ws.Range("N31").Value = -7495.4443
$ws.Range("H34").Value = 4443.5557
$ws.Range("I34").Value = 1981.6666
$ws.Range("J34").Value = 6905.4443
$ws.Range("K34").Value = 1981.6666
$ws.Range("L34").Value = 6905.4443
$ws.Range("M34").Value = -1779.6666
$ws.Range("N34").Value = -7309.4443
$ws.Range("H58").Value = 6168.625
$ws.Range("I58").Value = 4739.6
$ws.Range("J58").Value = 6818.1816
$ws.Range("K58").Value = 4739.6
$ws.Range("L58").Value = 6818.1816
$ws.Range("M58").Value = -4536.6
$ws.Range("N58").Value = -7224.1816
$ws.Range("H132").Value = 2926.353
$ws.Range("I132").Value = 2661.1936
$ws.Range("K132").Value = 7983.5808
$ws.Range("M132").Value = -5453.5808
$ws.Range("H134").Value = 3897.7964
$ws.Range("I134").Value = 2649.175
$ws.Range("K134").Value = 7947.525000000001
$ws.Range("M134").Value = -5412.525000000001
$ws.Range("H136").Value = 6168.625
$ws.Range("I136").Value = 4739.6
$ws.Range("J136").Value = 6818.1816
$ws.Range("K136").Value = 14218.8
$ws.Range("L136").Value = 20454.5448
$ws.Range("M136").Value = -11668.8
$ws.Range("N136").Value = -25554.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8199.4
$ws.Range("I5").Value = 8199.4
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 8199.4
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -8087.4
$ws.Range("N5").ClearContents()
$ws.Range("H70").Value = 15190.077
$ws.Range("I70").Value = 14134.818
$ws.Range("K70").Value = 14134.818
$ws.Range("M70").Value = -13864.818
$ws.Range("H73").Value = 15190.077
$ws.Range("I73").Value = 14134.818
$ws.Range("K73").Value = 14134.818
$ws.Range("M73").Value = -13198.818
$ws.Range("H122").Value = 5757.231
$ws.Range("I122").Value = 2768.8
$ws.Range("K122").Value = 8306.400000000001
$ws.Range("M122").Value = -5856.400000000001
$ws.Range("H126").Value = 4651.8335
$ws.Range("I126").Value = 2884.6
$ws.Range("J126").Value = 5914.143
$ws.Range("K126").Value = 8653.799999999999
$ws.Range("L126").Value = 17742.429
$ws.Range("M126").Value = -6183.799999999999
$ws.Range("N126").Value = -22682.429
$ws.Range("H132").Value = 1776
$ws.Range("I132").Value = 1655.3334
$ws.Range("K132").Value = 4966.0002
$ws.Range("M132").Value = -2436.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 57999
$ws.Range("J36").Value = 57999
$ws.Range("L36").Value = 57999
$ws.Range("N36").Value = -59123
$ws.Range("H40").Value = 9935.429
$ws.Range("I40").Value = 13654.444
$ws.Range("K40").Value = 13654.444
$ws.Range("M40").Value = -13518.444
$ws.Range("H68").Value = 5155.9287
$ws.Range("I68").Value = 4139
$ws.Range("J68").Value = 7698.25
$ws.Range("K68").Value = 4139
$ws.Range("L68").Value = 7698.25
$ws.Range("M68").Value = -3390
$ws.Range("N68").Value = -9196.25
$ws.Range("H71").Value = 5155.9287
$ws.Range("I71").Value = 4139
$ws.Range("J71").Value = 7698.25
$ws.Range("K71").Value = 20695
$ws.Range("L71").Value = 38491.25
$ws.Range("M71").Value = -16951
$ws.Range("N71").Value = -45979.25
$ws.Range("H82").Value = 3477.125
$ws.Range("I82").Value = 2025.2727
$ws.Range("K82").Value = 2025.2727
$ws.Range("M82").Value = -1664.2727
$ws.Range("H85").Value = 3477.125
$ws.Range("I85").Value = 2025.2727
$ws.Range("K85").Value = 2025.2727
$ws.Range("M85").Value = -777.2727
$ws.Range("H93").Value = 5666.3335
$ws.Range("I93").Value = 4999.5
$ws.Range("K93").Value = 4999.5
$ws.Range("M93").Value = -3751.5
$ws.Range("H132").Value = 5292
$ws.Range("I132").Value = 2615
$ws.Range("J132").Value = 6630.5
$ws.Range("K132").Value = 7845
$ws.Range("L132").Value = 19891.5
$ws.Range("M132").Value = -5315
$ws.Range("N132").Value = -24951.5
$ws.Range("H136").Value = 5019.269
$ws.Range("I136").Value = 3080.0667
$ws.Range("K136").Value = 9240.2001
$ws.Range("M136").Value = -6690.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 1500
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = -858
$ws.Range("N12").Value = -2284
$ws.Range("H31").Value = 20000
$ws.Range("J31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("N31").Value = -20696
$ws.Range("H122").Value = 4088.2917
$ws.Range("I122").Value = 3418.8235
$ws.Range("K122").Value = 10256.4705
$ws.Range("M122").Value = -7806.470499999999
$ws.Range("H132").Value = 3008.1365
$ws.Range("I132").Value = 1751.6364
$ws.Range("J132").Value = 4264.636
$ws.Range("K132").Value = 5254.9092
$ws.Range("L132").Value = 12793.908
$ws.Range("M132").Value = -2724.9092
$ws.Range("N132").Value = -17853.908
$ws.Range("H136").Value = 3685.9092
$ws.Range("I136").Value = 2108.75
$ws.Range("K136").Value = 6326.25
$ws.Range("M136").Value = -3776.25
